$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the header-column style (same as existing A3:A11) to the newly created rows A12:A26 first,
# so the row-style matches what Excel would carry when a user fills the column down.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A12:A35").PasteSpecial(-4122) | Out-Null

# Populate index labels + high/low/avg_close for fiscal years 1989-2022
$ws.Cells.Item(2, 1).Value = "AAPL-FY-1989"
$ws.Cells.Item(3, 1).Value = "AAPL-FY-1990"
$ws.Cells.Item(3, 2).Value = 0.4497770071029663
$ws.Cells.Item(3, 3).Value = 0.2433039993047714
$ws.Cells.Item(3, 4).Value = 0.354091441702275
$ws.Cells.Item(4, 1).Value = "AAPL-FY-1991"
$ws.Cells.Item(4, 2).Value = 0.654017984867096
$ws.Cells.Item(4, 3).Value = 0.2165179997682571
$ws.Cells.Item(4, 4).Value = 0.4320348132815626
$ws.Cells.Item(5, 1).Value = "AAPL-FY-1992"
$ws.Cells.Item(5, 2).Value = 0.625
$ws.Cells.Item(5, 3).Value = 0.3705359995365143
$ws.Cells.Item(5, 4).Value = 0.4837598591808259
$ws.Cells.Item(6, 1).Value = "AAPL-FY-1993"
$ws.Cells.Item(6, 2).Value = 0.5825889706611633
$ws.Cells.Item(6, 3).Value = 0.2098210006952286
$ws.Cells.Item(6, 4).Value = 0.4218617693002045
$ws.Cells.Item(7, 1).Value = "AAPL-FY-1994"
$ws.Cells.Item(7, 2).Value = 0.34375
$ws.Cells.Item(7, 3).Value = 0.1964289993047714
$ws.Cells.Item(7, 4).Value = 0.2819552562924713
$ws.Cells.Item(8, 1).Value = "AAPL-FY-1995"
$ws.Cells.Item(8, 2).Value = 0.4475449919700623
$ws.Cells.Item(8, 3).Value = 0.2901790142059326
$ws.Cells.Item(8, 4).Value = 0.3673780051961777
$ws.Cells.Item(9, 1).Value = "AAPL-FY-1996"
$ws.Cells.Item(9, 2).Value = 0.3794640004634857
$ws.Cells.Item(9, 3).Value = 0.1428570002317429
$ws.Cells.Item(9, 4).Value = 0.2511022906927836
$ws.Cells.Item(10, 1).Value = "AAPL-FY-1997"
$ws.Cells.Item(10, 2).Value = 0.2639510035514832
$ws.Cells.Item(10, 3).Value = 0.1138390004634857
$ws.Cells.Item(10, 4).Value = 0.1755720084858319
$ws.Cells.Item(11, 1).Value = "AAPL-FY-1998"
$ws.Cells.Item(11, 2).Value = 0.390625
$ws.Cells.Item(11, 3).Value = 0.1138390004634857
$ws.Cells.Item(11, 4).Value = 0.2294782953136471
$ws.Cells.Item(12, 1).Value = "AAPL-FY-1999"
$ws.Cells.Item(12, 2).Value = 0.7154020071029663
$ws.Cells.Item(12, 3).Value = 0.2544640004634857
$ws.Cells.Item(12, 4).Value = 0.3933011308785469
$ws.Cells.Item(13, 1).Value = "AAPL-FY-2000"
$ws.Cells.Item(13, 2).Value = 1.342633962631226
$ws.Cells.Item(13, 3).Value = 0.453125
$ws.Cells.Item(13, 4).Value = 0.9245598661760411
$ws.Cells.Item(14, 1).Value = "AAPL-FY-2001"
$ws.Cells.Item(14, 2).Value = 0.4842860102653503
$ws.Cells.Item(14, 3).Value = 0.2433039993047714
$ws.Cells.Item(14, 4).Value = 0.3565608266635462
$ws.Cells.Item(15, 1).Value = "AAPL-FY-2002"
$ws.Cells.Item(15, 2).Value = 0.4673210084438324
$ws.Cells.Item(15, 3).Value = 0.2464289963245392
$ws.Cells.Item(15, 4).Value = 0.3620105781165727
$ws.Cells.Item(16, 1).Value = "AAPL-FY-2003"
$ws.Cells.Item(16, 2).Value = 0.4164290130138397
$ws.Cells.Item(16, 3).Value = 0.2271430045366287
$ws.Cells.Item(16, 4).Value = 0.299975082396986
$ws.Cells.Item(17, 1).Value = "AAPL-FY-2004"
$ws.Cells.Item(17, 2).Value = 0.6941069960594177
$ws.Cells.Item(17, 3).Value = 0.34375
$ws.Cells.Item(17, 4).Value = 0.4770635803937912
$ws.Cells.Item(18, 1).Value = "AAPL-FY-2005"
$ws.Cells.Item(18, 2).Value = 1.921785950660706
$ws.Cells.Item(18, 3).Value = 0.65767902135849
$ws.Cells.Item(18, 4).Value = 1.325498145014521
$ws.Cells.Item(19, 1).Value = "AAPL-FY-2006"
$ws.Cells.Item(19, 2).Value = 3.085714101791382
$ws.Cells.Item(19, 3).Value = 1.709643006324768
$ws.Cells.Item(19, 4).Value = 2.346522038802505
$ws.Cells.Item(20, 1).Value = "AAPL-FY-2007"
$ws.Cells.Item(20, 2).Value = 5.535714149475098
$ws.Cells.Item(20, 3).Value = 2.592856884002686
$ws.Cells.Item(20, 4).Value = 3.712232831001282
$ws.Cells.Item(21, 1).Value = "AAPL-FY-2008"
$ws.Cells.Item(21, 2).Value = 7.248570919036865
$ws.Cells.Item(21, 3).Value = 4.122857093811035
$ws.Cells.Item(21, 4).Value = 5.841691822644724
$ws.Cells.Item(22, 1).Value = "AAPL-FY-2009"
$ws.Cells.Item(22, 2).Value = 6.746428966522217
$ws.Cells.Item(22, 3).Value = 2.792856931686401
$ws.Cells.Item(22, 4).Value = 4.284210282017985
$ws.Cells.Item(23, 1).Value = "AAPL-FY-2010"
$ws.Cells.Item(23, 2).Value = 10.48321437835693
$ws.Cells.Item(23, 3).Value = 6.45357084274292
$ws.Cells.Item(23, 4).Value = 8.192667900328617
$ws.Cells.Item(24, 1).Value = "AAPL-FY-2011"
$ws.Cells.Item(24, 2).Value = 15.10214328765869
$ws.Cells.Item(24, 3).Value = 9.821429252624512
$ws.Cells.Item(24, 4).Value = 12.21202808713156
$ws.Cells.Item(25, 1).Value = "AAPL-FY-2012"
$ws.Cells.Item(25, 2).Value = 25.18107032775879
$ws.Cells.Item(25, 3).Value = 12.65142917633057
$ws.Cells.Item(25, 4).Value = 18.81095007807016
$ws.Cells.Item(26, 1).Value = "AAPL-FY-2013"
$ws.Cells.Item(26, 2).Value = 24.16964340209961
$ws.Cells.Item(26, 3).Value = 13.75357055664062
$ws.Cells.Item(26, 4).Value = 17.2969492142459
$ws.Cells.Item(27, 1).Value = "AAPL-FY-2014"
$ws.Cells.Item(27, 2).Value = 25.93499946594238
$ws.Cells.Item(27, 3).Value = 16.94321441650391
$ws.Cells.Item(27, 4).Value = 20.90607998760573
$ws.Cells.Item(28, 1).Value = "AAPL-FY-2015"
$ws.Cells.Item(28, 2).Value = 33.6349983215332
$ws.Cells.Item(28, 3).Value = 23
$ws.Cells.Item(28, 4).Value = 29.64994020955971
$ws.Cells.Item(29, 1).Value = "AAPL-FY-2016"
$ws.Cells.Item(29, 2).Value = 30.95499992370605
$ws.Cells.Item(29, 3).Value = 22.36750030517578
$ws.Cells.Item(29, 4).Value = 26.19310755938648
$ws.Cells.Item(30, 1).Value = "AAPL-FY-2017"
$ws.Cells.Item(30, 2).Value = 41.23500061035156
$ws.Cells.Item(30, 3).Value = 26.02000045776367
$ws.Cells.Item(30, 4).Value = 34.14772453904152
$ws.Cells.Item(31, 1).Value = "AAPL-FY-2018"
$ws.Cells.Item(31, 2).Value = 57.41749954223633
$ws.Cells.Item(31, 3).Value = 37.56000137329102
$ws.Cells.Item(31, 4).Value = 45.58395429816379
$ws.Cells.Item(32, 1).Value = "AAPL-FY-2019"
$ws.Cells.Item(32, 2).Value = 58.36750030517578
$ws.Cells.Item(32, 3).Value = 35.5
$ws.Cells.Item(32, 4).Value = 48.00395001220703
$ws.Cells.Item(33, 1).Value = "AAPL-FY-2020"
$ws.Cells.Item(33, 2).Value = 137.9799957275391
$ws.Cells.Item(33, 3).Value = 53.15250015258789
$ws.Cells.Item(33, 4).Value = 80.69004953999918
$ws.Cells.Item(34, 1).Value = "AAPL-FY-2021"
$ws.Cells.Item(34, 2).Value = 157.2599945068359
$ws.Cells.Item(34, 3).Value = 107.3199996948242
$ws.Cells.Item(34, 4).Value = 131.0207171193157
$ws.Cells.Item(35, 1).Value = "AAPL-FY-2022"
$ws.Cells.Item(35, 2).Value = 182.9400024414062
$ws.Cells.Item(35, 3).Value = 129.0399932861328
$ws.Cells.Item(35, 4).Value = 158.6180480713863

Write-Output "done"
